# Apply the perturbation described by the commit "Updated the rest of the test files"
# to the optimization_parameters sheet (and related view/selection bookkeeping).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# optimization_parameters sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1 used to repeat the "value" header label across C1:F1 -- drop the
# redundant cells so only A1 ("optimization_parameter") / B1 ("value") remain.
$ws.Range("C1:F1").ClearContents()

# The "Model" row becomes "production_function" (value stays "Sigmoid").
$ws.Range("A8").Value = "production_function"

# A new "L_curve" parameter row is inserted right after it, defaulting to 0.
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0
$ws.Range("B9").NumberFormat = "0.00E+00"

# The old "Deletion" row (now pushed down to row 17 by the insert above) is
# removed entirely.
$ws.Rows.Item(17).Delete()

# This sheet becomes the active / selected tab, with a new selection range.
$ws.Activate()
$ws.Range("C1:G2").Select()
